# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures ---
# VALOR MORA total
$ws.Range("E11").Value = 284700
# Cant. Periodos
$ws.Range("F13").Value = 5

# --- Re-point the last table row's formatting to the "closing border" style
#     that used to live on row 21 (the row about to be retired), then drop
#     the now-duplicate row 21 and let everything below ride up by one. ---
$ws.Range("B21:J21").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$ws.Rows("21").Delete()

# --- Refresh the period values for the account-statement detail rows ---
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"
$ws.Range("E19").Value = "2507"
$ws.Range("E20").Value = "2508"
